$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "3000-30000-888812-18162001-310530-00000-00000-0000-0000-0000"
$ws.Range("D14").Value = "3000-30000-888812-18162001-310540-00000-00000-0000-0000-0000"
$ws.Range("D15").Value = "3000-30000-888812-52191003-310530-00000-00000-0000-0000-0000"
$ws.Range("D16").Value = "3000-30000-888812-52191003-310540-00000-00000-0000-0000-0000"

[void]$ws.Range("D17").Select()
